$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CH-Kt")
$ws.Columns("B:B").Insert()
Write-Output "done"
